$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("R2").Value = $null
$ws.Range("R3").Value = 2021
$ws.Range("R4").Value = 202551
$ws.Range("R5").Value = 2.9794303052841493

$ws.Range("R2").Style = $ws.Range("Q2").Style
$ws.Range("R3").Style = $ws.Range("Q3").Style
$ws.Range("R4").Style = $ws.Range("Q4").Style
$ws.Range("R5").Style = $ws.Range("Q5").Style

$ws.Range("R2").Select()
